$d = $word.ActiveDocument
Write-Host "Tables count:" $d.Tables.Count
Write-Host "Paragraphs count:" $d.Paragraphs.Count
